$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "N22DCPT053"
$ws.Range("C5").Value = "Thị"
$ws.Range("D5").Value = "Nơ"
$ws.Range("E5").Value = "D22CQPT01-N"
$ws.Range("F5").Value = 2

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "N22DCPT032"
$ws.Range("C6").Value = "Văn"
$ws.Range("D6").Value = "Từng"
$ws.Range("E6").Value = "D22CNQ01-T"
$ws.Range("F6").Value = 10

# Row 7
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "N21DCPT02"
$ws.Range("C7").Value = "Đinh"
$ws.Range("D7").Value = "Oanh"
$ws.Range("E7").Value = "D21CPPT-M"
$ws.Range("F7").Value = 3

$ws.Range("I6").Select()
